$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column C (was 23.5703125 -> target 55.7109375 character-width units)
$ws.Columns("C").ColumnWidth = 54.8333

# Move the active selection to B5 (was A26) and let the view scroll
# back so topLeftCell resets to the sheet's natural top-left (A1),
# matching the removal of topLeftCell="A8" in the sheetView.
$ws.Range("B5").Select()
